# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" (fund-level holdings) positioned
#    right before the "总计" (summary) sheet.
# 2) Prepend a new summary row for "2022-Q1" at the top of the "总计"
#    sheet's data (existing rows shift down by one).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create + position the new "2022-Q1" sheet
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q1"

$zj = $wb.Worksheets.Item("总计")
$newSheet.Move($zj)

# sheet handles can go stale across a reorder -- re-fetch by name
$newSheet = $wb.Worksheets.Item("2022-Q1")
$q4 = $wb.Worksheets.Item("2021-Q4")

# Carry over the look of the "2021-Q4" sheet: header row (B1:H1) and the
# index column (A2:A8) share its style.
$q4.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2:A8").Copy()
$newSheet.Range("A2:A8").PasteSpecial(-4122)

# -- header row --
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# -- index column (A2:A8) --
$newSheet.Range("A2").Value = 0
$newSheet.Range("A3").Value = 1
$newSheet.Range("A4").Value = 2
$newSheet.Range("A5").Value = 3
$newSheet.Range("A6").Value = 4
$newSheet.Range("A7").Value = 5
$newSheet.Range("A8").Value = 6

# -- fund code / name / scale / position figures: stored as text, so
#    force text formatting first (keeps leading zeros in fund codes and
#    matches the source data which stores these as plain strings) --
$textRange = $newSheet.Range("B2:G8")
$textRange.NumberFormat = "@"

$newSheet.Range("B2").Value = "870009"
$newSheet.Range("C2").Value = "广发资管平衡精选一年持有混合A"
$newSheet.Range("D2").Value = "11.34"
$newSheet.Range("E2").Value = "94.29"
$newSheet.Range("F2").Value = "4.17"
$newSheet.Range("G2").Value = "0.4729"

$newSheet.Range("B3").Value = "872019"
$newSheet.Range("C3").Value = "广发资管平衡精选一年持有混合C"
$newSheet.Range("D3").Value = "1.54"
$newSheet.Range("E3").Value = "94.29"
$newSheet.Range("F3").Value = "4.17"
$newSheet.Range("G3").Value = "0.0642"

$newSheet.Range("B4").Value = "011431"
$newSheet.Range("C4").Value = "泰达宏利消费服务混合A"
$newSheet.Range("D4").Value = "1.61"
$newSheet.Range("E4").Value = "81.15"
$newSheet.Range("F4").Value = "3.71"
$newSheet.Range("G4").Value = "0.0597"

$newSheet.Range("B5").Value = "519959"
$newSheet.Range("C5").Value = "长信多利灵活配置混合"
$newSheet.Range("D5").Value = "1.45"
$newSheet.Range("E5").Value = "85.11"
$newSheet.Range("F5").Value = "4.00"
$newSheet.Range("G5").Value = "0.0580"

$newSheet.Range("B6").Value = "013488"
$newSheet.Range("C6").Value = "长信多利灵活配置混合D"
$newSheet.Range("D6").Value = "1.45"
$newSheet.Range("E6").Value = "85.11"
$newSheet.Range("F6").Value = "4.00"
$newSheet.Range("G6").Value = "0.0580"

$newSheet.Range("B7").Value = "519987"
$newSheet.Range("C7").Value = "长信恒利优势混合"
$newSheet.Range("D7").Value = "0.22"
$newSheet.Range("E7").Value = "82.39"
$newSheet.Range("F7").Value = "4.09"
$newSheet.Range("G7").Value = "0.0090"

$newSheet.Range("B8").Value = "011432"
$newSheet.Range("C8").Value = "泰达宏利消费服务混合C"
$newSheet.Range("D8").Value = "0.15"
$newSheet.Range("E8").Value = "81.15"
$newSheet.Range("F8").Value = "3.71"
$newSheet.Range("G8").Value = "0.0056"

# drop the temporary "@" format again -- ClearFormats resets the style
# to the default (unstyled) xf while the stored cell keeps its text type
$textRange.ClearFormats()

# -- 仓位排名 column: numeric --
$newSheet.Range("H2").Value = 9
$newSheet.Range("H3").Value = 9
$newSheet.Range("H4").Value = 3
$newSheet.Range("H5").Value = 8
$newSheet.Range("H6").Value = 8
$newSheet.Range("H7").Value = 8
$newSheet.Range("H8").Value = 3

# ---------------------------------------------------------------------
# Step 2: prepend the "2022-Q1" row to the "总计" summary sheet
# ---------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")
$zj.Rows(2).Insert()

# the inserted row inherits stray formatting from the row below -- clear
# it so B2:D2 end up with the same (unstyled) look as the rest of the
# data rows
$zj.Range("B2:D2").ClearFormats()

# A2 keeps the same "index column" style as the rows beneath it
$zj.Range("A3").Copy()
$zj.Range("A2").PasteSpecial(-4122)

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q1"
$zj.Range("C2").Value = 7
$zj.Range("D2").Value = 0.73

# renumber the existing rows that shifted down
$zj.Range("A3").Value = 1
$zj.Range("A4").Value = 2
